$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target data for rows 2-11 (columns A-G)
$data = @{
    2  = @(37378, "Lívia Fernandes", "Financeiro", "Viagem de negócios", 8, 45097, 2787.71)
    3  = @(47690, "Theo da Rosa", "Jurídico", "Problemas pessoais", 8, 45079, 7111)
    4  = @(40490, "Brenda Pinto", "P&D", "Viagem de negócios", 7, 45095, 2734.46)
    5  = @(76346, "Emanuelly da Mota", "Atendimento ao Cliente", "Outros", 8, 45084, 7616.2)
    6  = @(1293, "Sr. Davi Luiz Silva", "Engenharia", "Outros", 6, 45079, 10202.58)
    7  = @(91398, "Cauã Nogueira", "Marketing", "Problemas pessoais", 1, 45090, 6254.44)
    8  = @(92351, "Rafael Jesus", "Vendas", "Consulta médica", 3, 45104, 3078.63)
    9  = @(93715, "Davi Lucca Costela", "Operações", "Outros", 4, 45095, 9206.18)
    10 = @(11675, "Daniel Castro", "Engenharia", "Viagem de negócios", 7, 45082, 11961.83)
    11 = @(51298, "Carlos Eduardo Cunha", "Marketing", "Problemas pessoais", 7, 45104, 11213.33)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $ws.Cells.Item($row, 1).Value = $values[0]
    $ws.Cells.Item($row, 2).Value = $values[1]
    $ws.Cells.Item($row, 3).Value = $values[2]
    $ws.Cells.Item($row, 4).Value = $values[3]
    $ws.Cells.Item($row, 5).Value = $values[4]
    $ws.Cells.Item($row, 6).Value = $values[5]
    $ws.Cells.Item($row, 7).Value = $values[6]
}
